$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.052.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.30%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.918.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.58%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.88%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.25%  "

# Row 7
$ws.Range("E7").Value = "  -0.14%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.917.46"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.55%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.502"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.39%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.44%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.145"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.14%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.448"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.44%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000225"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.17%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.57"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.18%  "

# Row 15
$ws.Range("E15").Value = "  +0.06%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.403.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.65%  "

# Row 17
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.86"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.46%  "

# Row 18
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.050.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.31%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.917.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.85%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "426.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.35%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.69%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.673"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.41%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.40%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.91%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.77%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.94%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.77%  "

# Row 28
$ws.Range("E28").Value = "  +0.00%  "

# Row 29
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.00%  "

# Row 30
$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.33%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.10%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.92%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.91%  "

# Row 34
$ws.Range("E34").Value = "  -2.82%  "

# Row 35
$ws.Range("E35").Value = "  +1.19%  "

# Row 36
$ws.Range("E36").Value = "  -0.35%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.75%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.22%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.78"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.63%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.80%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.125"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.63%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.91%  "

# Row 43
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.288"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.92%  "

# Row 44
$ws.Range("B44").Value = "Arweave"
$ws.Range("C44").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.45%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "378.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.52%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0347"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.48%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.655.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.62%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.07%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.42%  "

# Row 51
$ws.Range("E51").Value = "  -0.34%  "
